$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: politeness_score (column B) was stored as text "3"; convert it to a real number
$ws.Cells.Item(26, 2).Value = 3

# Row 27: new annotation row appended below the existing data
$ws.Cells.Item(27, 1).Value = "Sunsi Wu"

# Keep the politeness_score in this new row as literal text "3" (matches the source data)
$ws.Cells.Item(27, 2).NumberFormat = "@"
$ws.Cells.Item(27, 2).Value = "3"
$ws.Cells.Item(27, 2).Style = "Normal"

$ws.Cells.Item(27, 3).Value = "无"
$ws.Cells.Item(27, 4).Value = "SUG"
$ws.Cells.Item(27, 5).Value = "MET"
$ws.Cells.Item(27, 6).Value = "7e3af7f2-ed39-457b-b159-8a754cc477a9"
$ws.Cells.Item(27, 7).Value = "HkwZSG-CZ_annotated.xlsx"
$ws.Cells.Item(27, 8).Value = "In general, computational wall time of MoS is actually sub-linear w.r.t. the number of mixture components."
